# Updates cryptos list Price (D) / Volume(1h) (E) columns with latest scraped values.
# Leading "'" forces text storage for column D so multi-dot "numbers" and trailing
# zeros (e.g. "28.90") are preserved exactly as plain text, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.699.70'
$ws.Range('D3').Value = '''1.617.19'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = '''212.67'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '''0.521'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = '''0.993'
$ws.Range('E7').Value = '  -0.42%  '
$ws.Range('D8').Value = '''28.90'
$ws.Range('E8').Value = '  +8.24%  '
$ws.Range('E9').Value = '  +2.74%  '
$ws.Range('D10').Value = '''0.0609'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('D11').Value = '''0.0909'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '''1.849.59'
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('D13').Value = '''1.620.72'
$ws.Range('E13').Value = '  +1.39%  '
$ws.Range('E14').Value = '  +6.05%  '
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('D16').Value = '''29.681.03'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').Value = '''8.90'
$ws.Range('E17').Value = '  +16.45%  '
$ws.Range('D18').Value = '''64.43'
$ws.Range('E18').Value = '  +1.62%  '
$ws.Range('D19').Value = '''244.23'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').Value = '''0.0₃0710'
$ws.Range('E20').Value = '  +2.74%  '
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').Value = '''4.11'
$ws.Range('E22').Value = '  +3.12%  '
$ws.Range('E23').Value = '  +5.39%  '
$ws.Range('D24').Value = '''2.11'
$ws.Range('E24').Value = '  +0.89%  '
$ws.Range('D25').Value = '''156.64'
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').Value = '''15.66'
$ws.Range('E26').Value = '  +2.18%  '
$ws.Range('E27').Value = '  +1.37%  '
$ws.Range('E28').Value = '  +3.08%  '
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').Value = '''0.0483'
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('D31').Value = '''3.32'
$ws.Range('E31').Value = '  +2.98%  '
$ws.Range('D32').Value = '''1.07'
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('D33').Value = '''3.21'
$ws.Range('E33').Value = '  +3.24%  '
$ws.Range('D34').Value = '''1.435.64'
$ws.Range('E34').Value = '  +1.34%  '
$ws.Range('E35').Value = '  +6.86%  '
$ws.Range('D36').Value = '''1.05'
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('D37').Value = '''2.88'
$ws.Range('E37').Value = '  +2.82%  '
$ws.Range('E38').Value = '  -0.81%  '
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('D40').Value = '''0.554'
$ws.Range('E40').Value = '  +3.40%  '
$ws.Range('E41').Value = '  +4.00%  '
$ws.Range('D42').Value = '''0.828'
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = '''69.50'
$ws.Range('E44').Value = '  +5.96%  '
$ws.Range('D45').Value = '''53.67'
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('E46').Value = '  -0.51%  '
$ws.Range('E47').Value = '  +19.84%  '
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('D49').Value = '''1.758.84'
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('D50').Value = '''87.83'
$ws.Range('E50').Value = '  +1.39%  '
$ws.Range('E51').Value = '  -0.56%  '
